# Implement data driven testing in login page
#
# Renames the header labels and replaces the sample username/password rows
# with data-driven-testing style values (a "valid" credential row and an
# "Invalid" variant row), giving the data rows the familiar
# variable/data-token look: Consolas font on a light blue fill, green text
# for the valid value and blue text for the rest.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:B1) -------------------------------------------------
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "password"

# --- Data rows (A2:B3) ---------------------------------------------------
$ws.Range("A2").Value = "QualityWay"
$ws.Range("B2").Value = "pwd_ds_algo@2"
$ws.Range("A3").Value = "QualityWayInvalid"
$ws.Range("B3").Value = "pwd_ds_algo@2Invalid"

# Light blue fill (FFE8F2FE) + Consolas font on every data cell, with the
# "valid" username (A2) in green (FF3F7F5F) and the rest in blue (FF2A00FF).
$ws.Range("A2").Interior.Color = 16708328
$ws.Range("A2").Font.Name = "Consolas"
$ws.Range("A2").Font.Color = 6258495

$ws.Range("B2").Interior.Color = 16708328
$ws.Range("B2").Font.Name = "Consolas"
$ws.Range("B2").Font.Color = 16711722

$ws.Range("A3").Interior.Color = 16708328
$ws.Range("A3").Font.Name = "Consolas"
$ws.Range("A3").Font.Color = 16711722

$ws.Range("B3").Interior.Color = 16708328
$ws.Range("B3").Font.Name = "Consolas"
$ws.Range("B3").Font.Color = 16711722
